$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 24, shifting existing rows 24:83 down to 25:84
$ws.Rows.Item(24).Insert()

# Populate the newly inserted row 24 with the new record
$ws.Cells.Item(24, 1).Value = 1
$ws.Cells.Item(24, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(24, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(24, 4).Value = 45037
$ws.Cells.Item(24, 5).Value = 15
$ws.Cells.Item(24, 6).Value = 100112027
$ws.Cells.Item(24, 7).Value = "Melón"
$ws.Cells.Item(24, 8).Value = "Tuna"
$ws.Cells.Item(24, 9).Value = "Primera"
$ws.Cells.Item(24, 10).Value = 40
$ws.Cells.Item(24, 11).Value = 15000
$ws.Cells.Item(24, 12).Value = 16000
$ws.Cells.Item(24, 13).Value = 15625
$ws.Cells.Item(24, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(24, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(24, 16).Value = 868
$ws.Cells.Item(24, 17).Value = 18
$ws.Cells.Item(24, 18).Value = "Hortaliza"
